$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for every data row (rows 2-221).
# Update all of them from 45182 (2023-09-13) to 45184 (2023-09-15).
$ws.Range("C2:C221").Value = 45184
